$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion text (A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = $ws1.Range("A1").Value()
$text = $text.Replace("1000 Bs = 4.66 = 18384.54 pesos", "1000 Bs = 4.75 = 18764.85 pesos")
$text = $text.Replace("18384.54 pesos = 4.61 = 913.7 Bs", "18764.85 pesos = 4.71 = 961.82 Bs")
$ws1.Range("A1").Value = $text

# --- tasas: update rate figures in N10/O10/N12/O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 210.5
$ws2.Range("O10").Value = 3950
$ws2.Range("N12").Value = 3980
$ws2.Range("O12").Value = 204
